$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B used to hold the numeric "n" values; it now holds the study
# (trial) name text, with header "study".
$ws.Range("B1").Value = "study"
$ws.Range("B1").Font.Bold = $false

$ws.Range("B2").Value = "CTSN Severe MR"
$ws.Range("B3").Value = "CTSN Moderate MR"
$ws.Range("B4").Value = "CTSN TR Trial"
$ws.Range("B5").Value = "CTSN AF Trial"

# Rename header of column A: "study" -> "author_year"
$ws.Range("A1").Value = "author_year"

# Make the plot/table more compact: widen column B to fit the new text.
$ws.Columns.Item(2).ColumnWidth = 23.42578125

$ws.Range("F14").Select()
